$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep their original text (inline-string) representation
# rather than being auto-coerced into numeric values by Excel (e.g. "1.00" -> 1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.718.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.99%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +14.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.10"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +14.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.42"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.660.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.310.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.79%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.713.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +22.86%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.64%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.51"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.51%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "44.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +19.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.82"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +9.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0803"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.82%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +13.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.92%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.72"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +20.72%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +16.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.52"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +14.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0306"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.869.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +15.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.44"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +20.14%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +11.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "74.68"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +14.52%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +13.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.09"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.36%  "
